$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("draculaV2")
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 230
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Columns.Item(1).Hidden = $false
$ws1.Columns.Item(2).Hidden = $false
$ws1.Columns.Item(1).ColumnWidth = 7.5546875
$ws1.Columns.Item(2).ColumnWidth = 14.109375
